$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 93
$ws.Range("E17").Value = 58
$ws.Range("F17").Value = 22
$ws.Range("H17").Value = 22
$ws.Range("E18").Value = 48
$ws.Range("E25").Value = 9
$ws.Range("E27").Value = 4
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("E37").Value = 21
$ws.Range("F37").Value = 9
$ws.Range("H37").Value = 9
$ws.Range("E38").Value = 36
$ws.Range("F38").Value = 7
$ws.Range("H38").Value = 7
$ws.Range("E39").Value = 11
$ws.Range("E43").Value = 11
$ws.Range("E46").Value = 12
$ws.Range("E48").Value = 10
$ws.Range("E61").Value = 13
$ws.Range("E70").Value = 16
$ws.Range("E76").Value = 25
